# "updated RTMF to take out mode shifting"
# RTMF-freight: HDVs row used to shift 90% to rail and 10% to ships;
# remove that mode shifting (set both back to 0).

$wb = $excel.ActiveWorkbook
$wsPassengers = $wb.Worksheets.Item("RTMF-passengers")
$wsFreight = $wb.Worksheets.Item("RTMF-freight")

# HDVs row (row 3): rail (E) and ships (F) shift fractions -> 0
$wsFreight.Range("E3").Value = 0
$wsFreight.Range("F3").Value = 0

# Match the author's resulting selection/active-sheet state
$wsPassengers.Range("I13").Select()
$wsFreight.Activate()
$wsFreight.Range("F4").Select()
